$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Swap the K2 / L2 header labels (LogProx Internal <-> LogProx Internal +
#    External)
# ---------------------------------------------------------------------------
$k2Old = $ws.Range("K2").Value2
$l2Old = $ws.Range("L2").Value2
$ws.Range("K2").Value = $l2Old
$ws.Range("L2").Value = $k2Old

# ---------------------------------------------------------------------------
# 2. Add the new row 34 (copy formatting from row 33, then set its values).
#    Doing this before touching G20/G21 keeps the shared-string insertion
#    order the same as the authored edit (99% string first, then 100%, then
#    80%).
# ---------------------------------------------------------------------------
$ws.Range("A33:T33").Copy() | Out-Null
$ws.Range("A34:T34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(34).RowHeight = 43.2

$ws.Range("A34").Value = "100% NFE"
$ws.Range("B34").Value = 0.001
$ws.Range("C34").Value = 0.99
$ws.Range("D34").Value = 0.99
$ws.Range("E34").Value = "RAREsim-Python"
$ws.Range("F34").Value = "Hapgen Haplotype"
$ws.Range("G34").Value = "Separately-`n--functional_bins 99%`n--synonymous_bins 99%"
$ws.Range("H34").Value = 0.05
$ws.Range("I34").Value = 0.06
$ws.Range("J34").Value = 0.04
$ws.Range("K34").Value = 0.06
$ws.Range("L34").Value = 0.06
$ws.Range("M34").Value = 0.03
$ws.Range("N34").Value = 0.04
$ws.Range("O34").Value = 0.03
$ws.Range("P34").Value = 0.03
$ws.Range("Q34").Value = 0.08
$ws.Range("R34").Value = 0.04
$ws.Range("S34").Value = 0
$ws.Range("T34").Value = 0.04

# ---------------------------------------------------------------------------
# 3. Update G20 and G21 with the new, more specific pipeline descriptions.
# ---------------------------------------------------------------------------
$ws.Range("G20").Value = "Separately-`n--functional_bins 100%`n--synonymous_bins 100%"
$ws.Range("G21").Value = "Separately-`n--functional_bins 80%`n--synonymous_bins 80%"

# ---------------------------------------------------------------------------
# 4. Update the view so the new row is visible (best effort - selection is
#    persisted; the scrolled top-left row follows the selection/freeze).
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
$ws.Range("F37").Select()
